$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the localized title strings in row 1 (A1/B1/C1) ---
$ws.Range("A1").Value = "16.5.1.1a ""Аткаруу бийлигинин мамлекеттик органдарындагы жана жергиликтүү өз алдынча башкаруу органдарындагы коррупциянын деңгээли жөнүндө жеке түшүнүк"" индекси"
$ws.Range("B1").Value = "16.5.1.1a Индекс ""Личное представление об уровне коррупции в государственных органах исполнительной власти и органах местного самоуправления''"
$ws.Range("C1").Value = "16.5.1.1a Index ""Personal views about the level of corruption in executive government authorities and local government''"

# --- Add the new 2020 column (I) ---
$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").Value = 2020

# Row 5 - bold header-style data row
$ws.Range("H5").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("I5").Value = 12.3
$ws.Range("I5").NumberFormat = "0.0"

# Row 6 - normal data row (style reused for rows 7-13)
$ws.Range("H6").Copy()
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("I6").Value = 40.299999999999997
$ws.Range("I6").NumberFormat = "0.0"

$ws.Range("I6").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("I7").Value = 36.2
$ws.Range("I8").PasteSpecial(-4122)
$ws.Range("I8").Value = 44.3
$ws.Range("I9").PasteSpecial(-4122)
$ws.Range("I9").Value = 36
$ws.Range("I10").PasteSpecial(-4122)
$ws.Range("I10").Value = 2.7
$ws.Range("I11").PasteSpecial(-4122)
$ws.Range("I11").Value = 32.9
$ws.Range("I12").PasteSpecial(-4122)
$ws.Range("I12").Value = 11.3
$ws.Range("I13").PasteSpecial(-4122)
$ws.Range("I13").Value = -18.2

# Row 14 - bottom bordered total row
$ws.Range("H14").Copy()
$ws.Range("I14").PasteSpecial(-4122)
$ws.Range("I14").Value = 33
$ws.Range("I14").NumberFormat = "0.0"

$excel.CutCopyMode = $false

# --- Update the active selection to match the authored state ---
$ws.Range("F16").Select()
